# Opencart_testData.xlsx edit: refresh the test data values on
# userRegistrationData and add a new (empty) "LoginData" worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("userRegistrationData")

# Row 2 - Chaitanya / Somawar
$ws.Range("C2").Value = "chaitanya234aa45@gmail.com"
$ws.Range("D2").Value = "'334534112"
$ws.Range("E2").Value = "test@322121"

# Row 3 - Madhurima / lastname01
$ws.Range("C3").Value = "madhu898jkhkj179@gmail.com"
$ws.Range("D3").Value = "'452452323233"
$ws.Range("E3").Value = "test@23412122"

# Row 4 - Vaibhav / lastname02
$ws.Range("C4").Value = "vaibhav9utiut6549732@gmail.com"
$ws.Range("D4").Value = "'9879378711"
$ws.Range("E4").Value = "test@12123"

# Row 5 - Nikhil / lastname03
$ws.Range("C5").Value = "nikhilq3hjkytu7653445@test.com"
$ws.Range("D5").Value = "'232387987"
$ws.Range("E5").Value = "test@12124"

# Add the new "LoginData" worksheet right after the existing sheet.
$newSheet = $wb.Worksheets.Add([System.Type]::Missing, $ws)
$newSheet.Name = "LoginData"

# Re-select userRegistrationData and move the selection cursor to C5,
# matching the saved state of the workbook.
$ws.Activate()
$ws.Range("C5").Select()
